# Add new question row (297. Serialize and Deserialize Binary Tree)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (26) down to the new row (27)
$ws.Range("A26:C26").Copy()
$ws.Range("A27:C27").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row's values
$ws.Cells.Item(27, 1).Value = 297
$ws.Cells.Item(27, 2).Value = "NA"
$ws.Cells.Item(27, 3).Value = "Serialize and Deserialize Binary Tree"

# Match the resulting selection left behind by the edit
$ws.Range("F30").Select() | Out-Null
